$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$shp = $m.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
Write-Output ("len=" + $tr.Length)
$c = $tr.Characters(1, 3)
Write-Output ("partial=" + $c.Text)
$c.Text = "6/5"
Write-Output ("after=" + $shp.TextFrame.TextRange.Text)
